# Rebuild Finance templates with correct industry content (Banking)
$wb = $excel.ActiveWorkbook

$wsOverview   = $wb.Worksheets.Item("Change Management Overview")
$wsImpact     = $wb.Worksheets.Item("Change Impact Assessment")
$wsActivities = $wb.Worksheets.Item("Change Activities")

# --- Sheet: Change Management Overview ---
$wsOverview.Range("A2").Value = "Banking Implementation Project"
$wsOverview.Range("B6").Value = "Enterprise Banking Implementation"

# touch row 13 so it materializes as a blank row between the project info
# block (ends row 12) and the objectives header (row 14)
$wsOverview.Rows.Item(13).OutlineLevel = 0

$wsOverview.Range("A15").Value = "1. Achieve 95% user adoption of new Banking systems within 6 months of go-live"
$wsOverview.Range("A17").Value = "3. Build organizational capability and confidence in Banking technologies"
$wsOverview.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for Banking transformation"

# touch row 21 so it materializes as a blank row between the objectives
# block (ends row 20) and the strategy header (row 22)
$wsOverview.Rows.Item(21).OutlineLevel = 0

# --- Sheet: Change Impact Assessment ---
# touch row 2 so it materializes as a blank row between the title (row 1)
# and the table header (row 3)
$wsImpact.Rows.Item(2).OutlineLevel = 0

$wsImpact.Range("G4").Value = "Banking automation"

# --- Sheet: Change Activities ---
# touch row 2 so it materializes as a blank row between the title (row 1)
# and the table header (row 3)
$wsActivities.Rows.Item(2).OutlineLevel = 0
